$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 879
$ws.Cells.Item(3, 6).Value = 1458
$ws.Cells.Item(4, 6).Value = 1112
$ws.Cells.Item(5, 6).Value = 529
$ws.Cells.Item(7, 6).Value = 3
$ws.Cells.Item(8, 6).Value = 680
$ws.Cells.Item(11, 6).Value = 92
$ws.Cells.Item(12, 6).Value = 220
$ws.Cells.Item(13, 6).Value = 160
$ws.Cells.Item(14, 6).Value = 2689
$ws.Cells.Item(15, 6).Value = 2
$ws.Cells.Item(16, 6).Value = 6
$ws.Cells.Item(19, 6).Value = 505
$ws.Cells.Item(21, 6).Value = 410
$ws.Cells.Item(24, 6).Value = 668
$ws.Cells.Item(26, 6).Value = 248
$ws.Cells.Item(27, 6).Value = 968
$ws.Cells.Item(29, 6).Value = 1588
$ws.Cells.Item(30, 6).Value = 323

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(5, 6).Value = 232
$ws.Cells.Item(7, 6).Value = 222
$ws.Cells.Item(8, 6).Value = 285
$ws.Cells.Item(11, 6).Value = 29
$ws.Cells.Item(12, 6).Value = 132

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(3, 6).Value = 46

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 879
$ws.Cells.Item(4, 6).Value = 1458
$ws.Cells.Item(5, 6).Value = 1112
$ws.Cells.Item(8, 6).Value = 46
$ws.Cells.Item(9, 6).Value = 529
$ws.Cells.Item(11, 6).Value = 3
$ws.Cells.Item(12, 6).Value = 680
$ws.Cells.Item(16, 6).Value = 92
$ws.Cells.Item(17, 6).Value = 220
$ws.Cells.Item(18, 6).Value = 160
$ws.Cells.Item(19, 6).Value = 2689
$ws.Cells.Item(20, 6).Value = 2
$ws.Cells.Item(21, 6).Value = 6
$ws.Cells.Item(22, 6).Value = 232
$ws.Cells.Item(25, 6).Value = 505
$ws.Cells.Item(27, 6).Value = 410
$ws.Cells.Item(31, 6).Value = 222
$ws.Cells.Item(32, 6).Value = 285
$ws.Cells.Item(34, 6).Value = 668
$ws.Cells.Item(36, 6).Value = 29
$ws.Cells.Item(37, 6).Value = 132
$ws.Cells.Item(38, 6).Value = 132
$ws.Cells.Item(40, 6).Value = 248
$ws.Cells.Item(41, 6).Value = 968
$ws.Cells.Item(43, 6).Value = 1588
$ws.Cells.Item(44, 6).Value = 323
